$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A36").Value = "columnsProjectionsHighEdu"
$ws1.Range("B36").Value = 2

$ws1.Range("A37").Value = "columnsProjectionsLowEdu"
$ws1.Range("B37").Value = 2

$ws1.Range("A38").Value = "columnsStudentShareProjections"
$ws1.Range("B38").Value = 40

$ws1.Range("A39").Value = "columnsEmploymentAlignment"
$ws1.Range("B39").Value = 40

$ws1.Range("A40").Value = "columnsFertilityProjectionsByYear"
$ws1.Range("B40").Value = 90

$ws1.Range("A41").Value = "columnsCoefficientMapRMSE"
$ws1.Range("B41").Value = 1

$ws1.Range("A42").Value = "columnsMortalityProbabilityByGenderAgeYear"
$ws1.Range("B42").Value = 111

$ws1.Range("A43").Value = "columnsPopulationProjections"
$ws1.Range("B43").NumberFormat = "@"
$ws1.Range("B43").Value = "50"
$ws1.Range("B43").NumberFormat = "0"

$ws1.Activate()
$ws1.Range("A36:B43").Select()
